$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Cabernet Franc", 1300, "Bodega Cinco", "Famoso por sus vinos tintos de alta calidad", "Argentina", "Cabernet FrancSauvignon Blanc es una variedad de uva blanca conocida por su frescura y sus aromas herbáceos y cítricos.Bonarda es una variedad de uva tinta que se utiliza en la producción de vinos tintos suaves y afrutados, con sabores a frutas negras y especias.Sangiovese es una variedad de uva tinta que se asocia principalmente con los vinos italianos, conocidos por su acidez y sabor a frutas rojas.", 9.9),
    @("Malbec", 1200, "Bodega Cinco", "Famoso por sus vinos tintos de alta calidad", "Argentina", "MalbecMalbec", 9.199999999999999),
    @("Torrontés Clásico", 1150, "Bodega Dos", "Reconocida por sus Malbecs", "Argentina", "Torrontés ClásicoMourvèdre es una variedad de uva tinta que se utiliza en la producción de vinos tintos robustos y especiados.Cabernet Sauvignon es una variedad de uva tinta ampliamente reconocida por su presencia en los vinos tintos de Bordeaux.Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales.", 9),
    @("Merlot Reserva", 1250, "Bodega Tres", "Región importante de San Juan", "Argentina", "Merlot ReservaMerlot Reserva", 8.9),
    @("Riesling", 1350, "Bodega Ocho", "Famoso por sus vinos tintos de alta calidad", "Argentina", "RieslingRiesling", 8.6),
    @("Chardonnay", 1350, "Bodega Cuatro", "Famoso por sus vinos tintos de alta calidad", "Argentina", "Sémillon es una variedad de uva blanca que se utiliza en la producción de vinos blancos secos, dulces y también vinos de postre.Pinot Grigio es una variedad de uva blanca que produce vinos blancos ligeros y refrescantes, con notas cítricas y florales.Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales.", 8.4),
    @("Zinfandel", 1100, "Bodega Ocho", "Famoso por sus vinos tintos de alta calidad", "Argentina", "ZinfandelZinfandel", 8.324999999999999),
    @("Cabernet", 1300, "Bodega Seis", "Reconocida por sus Malbecs", "Argentina", "CabernetMerlot es una variedad de uva tinta que se caracteriza por su suavidad y sabor frutal en los vinos.Pinot Noir es una variedad de uva tinta que se asocia con vinos ligeros, elegantes y afrutados.", 8.300000000000001),
    @("Pinot Noir", 1500, "Bodega Tres", "Región importante de San Juan", "Argentina", "Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales.Tannat es una variedad de uva tinta que se asocia principalmente con los vinos de Uruguay, conocidos por su estructura tánica y sabor intenso.", 8.25),
    @("Torrontés", 1150, "Bodega Ocho", "Famoso por sus vinos tintos de alta calidad", "Argentina", "TorrontésRiesling es una variedad de uva blanca que puede producir desde vinos secos y refrescantes hasta vinos dulces y aromáticos.", 8.199999999999999),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}
